$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "Data"
$wsMeta = $wb.Worksheets.Item(2)
$wsMeta.Name = "Metadata"

# --- Sheet "Data": years now listed descending (2019 -> 2007), values follow the year ---
$dataRows = @(
    @("2019", 4.6, 4.5),
    @("2018", 4.7, 4.6),
    @("2017", 4.7, 4.6),
    @("2016", 5.1, 4.9),
    @("2015", 5.3, 4.9),
    @("2014", 5.1, 4.9),
    @("2013", 4.9, 4.7),
    @("2012", 5.8, 5.3),
    @("2011", 5.2, 4.8),
    @("2010", 6.1, 5.8),
    @("2009", 5.7, 5.1),
    @("2008", 5.8, 5.4),
    @("2007", 5.8, 5.2)
)

$r = 2
foreach ($row in $dataRows) {
    $cellA = $wsData.Cells.Item($r, 1)
    # Force the year label to stay text (it would otherwise be auto-detected
    # as a number) then drop the temporary "@" number format again so no
    # stray style is left behind.
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.ClearFormats()

    $wsData.Cells.Item($r, 2).Value = $row[1]
    $wsData.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Sheet "Metadata": reorder + rename labels to lowercase, add observaciones + footer rows ---
$wsMeta.Cells.Item(2, 1).Value = "nomindicador"
$wsMeta.Cells.Item(2, 2).Value = "Porcentaje de personas que viven en asentamientos"

$wsMeta.Cells.Item(3, 1).Value = "derecho"
$wsMeta.Cells.Item(3, 2).Value = "Vivienda"

$wsMeta.Cells.Item(4, 1).Value = "conindicador"
$wsMeta.Cells.Item(4, 2).Value = "Asentamientos"

$wsMeta.Cells.Item(5, 1).Value = "tipoind"
$wsMeta.Cells.Item(5, 2).Value = "Resultados"

$wsMeta.Cells.Item(6, 1).Value = "definicion"
$wsMeta.Cells.Item(6, 2).Value = "El indicador mide el porcentaje de personas en viviendas ubicadas en asentamiento irregular."

$wsMeta.Cells.Item(7, 1).Value = "calculo"
$wsMeta.Cells.Item(7, 2).Value = "Para cada año calcular: (Cantidad de personas que residen en viviendas ubicadas en asentamiento irregular / Cantidad total de personas en viviendas particulares)*100"

$wsMeta.Cells.Item(8, 1).Value = "observaciones"
$wsMeta.Cells.Item(8, 2).Value = "Sin observaciones"

$wsMeta.Cells.Item(9, 1).Value = "cita"
$wsMeta.Cells.Item(9, 2).Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"

$wsMeta.Cells.Item(10, 1).Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$wsMeta.Cells.Item(10, 2).Value = " "
